$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "b24d361c-eb46-42a9-9c4a-10ab55ca990d"
$ws.Range("B15").Value = "plywood"
$ws.Range("C15").Value = "p002"
$ws.Range("D15").Value = "2026-01-19T21:12:07.822Z"
$ws.Range("E15").Value = -30
$ws.Range("F15").Value = -10

$ws.Range("A16").Value = "3c895075-7194-4a91-b14b-2d30f67e4662"
$ws.Range("B16").Value = "hardware"
$ws.Range("C16").Value = "Hardware 1mm"
$ws.Range("D16").Value = "2026-01-19T21:19:04.819Z"
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 180
